# Update "Note/comments" text for TEND, IADAPTIME, STPMIN and STPMAX rows,
# adjust alignment / wrap-text formatting for the affected rows, tweak row
# heights and refresh the frozen-pane view position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Note/comments (column D) text -----------------------------
$ws.Range("D6").Value = "End time of the simulation"
$ws.Range("D7").Value = "Flag for time adaptivity: 0 = no adaptivity (tstep=tstepmin); 1 = adaptivity based on variation of the whole thermal-hydraulic solution; 2 = adaptivity based only on temperature variation of the thermal-hydraulic solution; -1 = adaptive time step from user defined input file (to be implemented); -2 = adaptive time step from user defined function."
$ws.Range("D8").Value = "Minimum time step for the thermal-hydraulic loop. Used if flag IADAPTIME = 0. Lower bound of the adaptive time step if IADAPTIME = 1 or IADAPTIME = 2"
$ws.Range("D9").Value = "Maximum time step for the thermal-hydraulic loop. It is the Upper bound of the adaptive time step if IADAPTIME = 1 or IADAPTIME = 2"

# --- Formatting: rows 7-9 get vertically-centered + wrapped text ----------
$rows789 = $ws.Range("A7:D9")
$rows789.VerticalAlignment = -4108   # xlCenter
$ws.Range("A7:A9").WrapText = $false
$ws.Range("B7:C9").HorizontalAlignment = -4108  # xlCenter
$ws.Range("D7:D9").WrapText = $true
$ws.Range("D10").WrapText = $true
$ws.Range("D10").VerticalAlignment = -4108

# --- Row heights -------------------------------------------------------
$ws.Rows("7").RowHeight = 58
$ws.Rows("8").RowHeight = 29
$ws.Rows("9").RowHeight = 29

# --- Frozen pane / selection bookkeeping -----------------------------
$ws.Range("D10").Select()
